$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -186
$ws.Range("N4").ClearContents()
$ws.Range("H28").Value = 1861.7778
$ws.Range("I28").Value = 2539.8
$ws.Range("J28").Value = 1014.25
$ws.Range("K28").Value = 2539.8
$ws.Range("L28").Value = 1014.25
$ws.Range("M28").Value = -2054.8
$ws.Range("N28").Value = -1984.25
$ws.Range("H40").Value = 3273.7
$ws.Range("I40").Value = 3342.5
$ws.Range("J40").Value = 2998.5
$ws.Range("K40").Value = 3342.5
$ws.Range("L40").Value = 2998.5
$ws.Range("M40").Value = -3167.5
$ws.Range("N40").Value = -3348.5
$ws.Range("H74").Value = 125010000
$ws.Range("I74").Value = 166673330
$ws.Range("K74").Value = 166673330
$ws.Range("M74").Value = -166672394
$ws.Range("H77").Value = 125010000
$ws.Range("I77").Value = 166673330
$ws.Range("K77").Value = 833366650
$ws.Range("M77").Value = -833361970
$ws.Range("H96").Value = 1282.1666
$ws.Range("I96").Value = 1166
$ws.Range("J96").Value = 1398.3334
$ws.Range("K96").Value = 3498
$ws.Range("L96").Value = 4195.0002
$ws.Range("M96").Value = -2125
$ws.Range("N96").Value = -6941.0002
$ws.Range("H100").Value = 3352.5715
$ws.Range("I100").Value = 1833
$ws.Range("J100").Value = 4492.25
$ws.Range("K100").Value = 1833
$ws.Range("L100").Value = 4492.25
$ws.Range("M100").Value = -1292
$ws.Range("N100").Value = -5574.25
$ws.Range("H111").Value = 41668984
$ws.Range("I111").Value = 62500500
$ws.Range("J111").Value = 5950
$ws.Range("K111").Value = 187501500
$ws.Range("L111").Value = 17850
$ws.Range("M111").Value = -187498433
$ws.Range("N111").Value = -23984
$ws.Range("H116").Value = 41673284
$ws.Range("I116").Value = 83337896
$ws.Range("J116").Value = 8670.666999999999
$ws.Range("K116").Value = 83337896
$ws.Range("L116").Value = 8670.666999999999
$ws.Range("M116").Value = -83334454
$ws.Range("N116").Value = -15554.667
$ws.Range("H132").Value = 1259.909
$ws.Range("I132").Value = 1259.909
$ws.Range("K132").Value = 3779.727
$ws.Range("M132").Value = -1249.727
$ws.Range("H137").Value = 3313.6667
$ws.Range("I137").Value = 3090.3333
$ws.Range("K137").Value = 9270.999899999999
$ws.Range("M137").Value = -6720.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1670614.1
$ws.Range("I32").Value = 1740197.4
$ws.Range("K32").Value = 1740197.4
$ws.Range("M32").Value = -1739910.4
$ws.Range("H45").Value = 5373.6875
$ws.Range("I45").Value = 2703.6
$ws.Range("J45").Value = 9823.833000000001
$ws.Range("K45").Value = 2703.6
$ws.Range("L45").Value = 9823.833000000001
$ws.Range("M45").Value = -2326.6
$ws.Range("N45").Value = -10577.833
$ws.Range("H74").Value = 57459.266
$ws.Range("I74").Value = 109389.266
$ws.Range("J74").Value = 5529.2666
$ws.Range("K74").Value = 109389.266
$ws.Range("L74").Value = 5529.2666
$ws.Range("M74").Value = -108515.266
$ws.Range("N74").Value = -7277.2666
$ws.Range("H77").Value = 57459.266
$ws.Range("I77").Value = 109389.266
$ws.Range("J77").Value = 5529.2666
$ws.Range("K77").Value = 546946.3300000001
$ws.Range("L77").Value = 27646.333
$ws.Range("M77").Value = -542578.3300000001
$ws.Range("N77").Value = -36382.333
$ws.Range("H97").Value = 27831116
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 27831116
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 27831116
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -27832108
$ws.Range("H102").Value = 1349.4615
$ws.Range("I102").Value = 1340.2727
$ws.Range("K102").Value = 1340.2727
$ws.Range("M102").Value = 281.7273
$ws.Range("H122").Value = 22665.455
$ws.Range("J122").Value = 9250
$ws.Range("L122").Value = 27750
$ws.Range("N122").Value = -32650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12822677
$ws.Range("I20").Value = 33336106
$ws.Range("J20").Value = 1783.75
$ws.Range("K20").Value = 33336106
$ws.Range("L20").Value = 1783.75
$ws.Range("M20").Value = -33335859
$ws.Range("N20").Value = -2277.75
$ws.Range("H86").Value = 45504456
$ws.Range("I86").Value = 86830.914
$ws.Range("J86").Value = 100005600
$ws.Range("K86").Value = 86830.914
$ws.Range("L86").Value = 100005600
$ws.Range("M86").Value = -85707.914
$ws.Range("N86").Value = -100007846
$ws.Range("H89").Value = 45504456
$ws.Range("I89").Value = 86830.914
$ws.Range("J89").Value = 100005600
$ws.Range("K89").Value = 434154.57
$ws.Range("L89").Value = 500028000
$ws.Range("M89").Value = -428538.57
$ws.Range("N89").Value = -500039232
$ws.Range("H94").Value = 1300.8695
$ws.Range("I94").Value = 479.73334
$ws.Range("K94").Value = 479.73334
$ws.Range("M94").Value = -28.73334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12826743
$ws.Range("I58").Value = 29413330
$ws.Range("K58").Value = 29413330
$ws.Range("M58").Value = -29413127
$ws.Range("H62").Value = 5850.364
$ws.Range("I62").Value = 6877.8335
$ws.Range("J62").Value = 4617.4
$ws.Range("K62").Value = 6877.8335
$ws.Range("L62").Value = 4617.4
$ws.Range("M62").Value = -6253.8335
$ws.Range("N62").Value = -5865.4
$ws.Range("H65").Value = 5850.364
$ws.Range("I65").Value = 6877.8335
$ws.Range("J65").Value = 4617.4
$ws.Range("K65").Value = 34389.1675
$ws.Range("L65").Value = 23087
$ws.Range("M65").Value = -31269.1675
$ws.Range("N65").Value = -29327
$ws.Range("H105").Value = 17867266
$ws.Range("I105").Value = 71428570
$ws.Range("K105").Value = 71428570
$ws.Range("M105").Value = -71426823
$ws.Range("H107").Value = 2342.8635
$ws.Range("I107").Value = 1631.9231
$ws.Range("K107").Value = 1631.9231
$ws.Range("M107").Value = 288.0769
$ws.Range("H134").Value = 6310.0293
$ws.Range("I134").Value = 1609.7142
$ws.Range("J134").Value = 9600.25
$ws.Range("K134").Value = 4829.142599999999
$ws.Range("L134").Value = 28800.75
$ws.Range("M134").Value = -2294.142599999999
$ws.Range("N134").Value = -33870.75
$ws.Range("H136").Value = 12826743
$ws.Range("I136").Value = 29413330
$ws.Range("K136").Value = 88239990
$ws.Range("M136").Value = -88237440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1438941.1
$ws.Range("I11").Value = 1917556.5
$ws.Range("K11").Value = 5752669.5
$ws.Range("M11").Value = -5752529.5
$ws.Range("H102").Value = 13000
$ws.Range("J102").Value = 13000
$ws.Range("L102").Value = 39000
$ws.Range("N102").Value = -43868

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 221.30302
$ws.Range("I2").Value = 192.7037
$ws.Range("K2").Value = 192.7037
$ws.Range("M2").Value = -79.7037
$ws.Range("H70").Value = 8482.5
$ws.Range("I70").Value = 5968
$ws.Range("K70").Value = 5968
$ws.Range("M70").Value = -5698
$ws.Range("H73").Value = 8482.5
$ws.Range("I73").Value = 5968
$ws.Range("K73").Value = 5968
$ws.Range("M73").Value = -5032
$ws.Range("H80").Value = 5498.5
$ws.Range("I80").Value = 4500
$ws.Range("K80").Value = 4500
$ws.Range("M80").Value = -3502
$ws.Range("H83").Value = 5498.5
$ws.Range("I83").Value = 4500
$ws.Range("K83").Value = 22500
$ws.Range("M83").Value = -17508
$ws.Range("H122").Value = 2133094.8
$ws.Range("I122").Value = 3151131.8
$ws.Range("J122").Value = 4471.909
$ws.Range("K122").Value = 9453395.399999999
$ws.Range("L122").Value = 13415.727
$ws.Range("M122").Value = -9450945.399999999
$ws.Range("N122").Value = -18315.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 35714660
$ws.Range("I55").Value = 76923150
$ws.Range("J55").Value = 629.3333
$ws.Range("K55").Value = 76923150
$ws.Range("L55").Value = 629.3333
$ws.Range("M55").Value = -76922977
$ws.Range("N55").Value = -975.3333
$ws.Range("H61").Value = 5362.72
$ws.Range("I61").Value = 3355.2144
$ws.Range("J61").Value = 7917.727
$ws.Range("K61").Value = 3355.2144
$ws.Range("L61").Value = 7917.727
$ws.Range("M61").Value = -3153.2144
$ws.Range("N61").Value = -8321.726999999999
$ws.Range("H100").Value = 4641.6665
$ws.Range("I100").Value = 3639.2
$ws.Range("K100").Value = 3639.2
$ws.Range("M100").Value = -3098.2
$ws.Range("H113").Value = 5362.72
$ws.Range("I113").Value = 3355.2144
$ws.Range("J113").Value = 7917.727
$ws.Range("K113").Value = 3355.2144
$ws.Range("L113").Value = 7917.727
$ws.Range("M113").Value = -1185.2144
$ws.Range("N113").Value = -12257.727
$ws.Range("H132").Value = 11370373
$ws.Range("I132").Value = 27780378
$ws.Range("J132").Value = 9599.962
$ws.Range("K132").Value = 83341134
$ws.Range("L132").Value = 28799.886
$ws.Range("M132").Value = -83338604
$ws.Range("N132").Value = -33859.886

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8665
$ws.Range("I62").Value = 8998
$ws.Range("J62").Value = 7999
$ws.Range("K62").Value = 8998
$ws.Range("L62").Value = 7999
$ws.Range("M62").Value = -8374
$ws.Range("N62").Value = -9247
$ws.Range("H65").Value = 8665
$ws.Range("I65").Value = 8998
$ws.Range("J65").Value = 7999
$ws.Range("K65").Value = 44990
$ws.Range("L65").Value = 39995
$ws.Range("M65").Value = -41870
$ws.Range("N65").Value = -46235
$ws.Range("H81").Value = 16673772
$ws.Range("I81").Value = 2525.9
$ws.Range("K81").Value = 5051.8
$ws.Range("M81").Value = -3990.8
$ws.Range("H84").Value = 16673772
$ws.Range("I84").Value = 2525.9
$ws.Range("K84").Value = 25259
$ws.Range("M84").Value = -19955
$ws.Range("H107").Value = 586.4783
$ws.Range("I107").Value = 389
$ws.Range("K107").Value = 1167
$ws.Range("M107").Value = 753
$ws.Range("H122").Value = 176521.12
$ws.Range("I122").Value = 251576.62
$ws.Range("J122").Value = 4965.7144
$ws.Range("K122").Value = 754729.86
$ws.Range("L122").Value = 14897.1432
$ws.Range("M122").Value = -752279.86
$ws.Range("N122").Value = -19797.1432
